$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name and link swaps (rows 25-28 and 50-51)
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('B28').Value = 'Dai'
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

# Price (column D) updates
$ws.Range('D2').Value = '62.359.80'
$ws.Range('D3').Value = '3.010.43'
$ws.Range('D9').Value = '3.012.01'
$ws.Range('D16').Value = '3.500.45'
$ws.Range('D18').Value = '62.270.96'
$ws.Range('D19').Value = '3.007.13'
$ws.Range('D36').Value = '0.0₃0798'
$ws.Range('D46').Value = '2.728.46'

# Numeric-looking prices need to be forced to text to match source formatting
$priceCells = @('D4','D5','D6','D8','D10','D11','D12','D13','D14','D17','D20','D21','D22','D23','D24','D25','D26','D27','D28','D29','D31','D32','D33','D34','D37','D38','D39','D40','D41','D42','D43','D44','D47','D48','D50','D51')
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '583.37'
$ws.Range('D6').Value = '146.36'
$ws.Range('D8').Value = '0.527'
$ws.Range('D10').Value = '0.149'
$ws.Range('D11').Value = '5.82'
$ws.Range('D12').Value = '0.465'
$ws.Range('D13').Value = '0.0000229'
$ws.Range('D14').Value = '34.57'
$ws.Range('D17').Value = '7.16'
$ws.Range('D20').Value = '461.91'
$ws.Range('D21').Value = '14.04'
$ws.Range('D22').Value = '0.690'
$ws.Range('D23').Value = '7.48'
$ws.Range('D24').Value = '81.72'
$ws.Range('D25').Value = '2.22'
$ws.Range('D26').Value = '12.33'
$ws.Range('D27').Value = '10.03'
$ws.Range('D28').Value = '1.00'
$ws.Range('D29').Value = '0.999'
$ws.Range('D31').Value = '7.03'
$ws.Range('D32').Value = '28.46'
$ws.Range('D33').Value = '2.09'
$ws.Range('D34').Value = '0.109'
$ws.Range('D37').Value = '5.78'
$ws.Range('D38').Value = '2.11'
$ws.Range('D39').Value = '9.31'
$ws.Range('D40').Value = '50.35'
$ws.Range('D41').Value = '2.88'
$ws.Range('D42').Value = '0.115'
$ws.Range('D43').Value = '393.35'
$ws.Range('D44').Value = '0.0358'
$ws.Range('D47').Value = '36.73'
$ws.Range('D48').Value = '128.55'
$ws.Range('D50').Value = '2.22'
$ws.Range('D51').Value = '0.109'
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# Volume(1h) (column E) updates
$ws.Range('E2').Value = '  -2.23%  '
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('E5').Value = '  -1.91%  '
$ws.Range('E6').Value = '  -5.42%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -2.60%  '
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('E10').Value = '  -5.02%  '
$ws.Range('E11').Value = '  -1.49%  '
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('E13').Value = '  -3.63%  '
$ws.Range('E14').Value = '  -6.12%  '
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('E18').Value = '  -2.20%  '
$ws.Range('E19').Value = '  -2.14%  '
$ws.Range('E20').Value = '  -4.58%  '
$ws.Range('E21').Value = '  -3.25%  '
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('E25').Value = '  -10.06%  '
$ws.Range('E26').Value = '  -4.63%  '
$ws.Range('E27').Value = '  -5.74%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  -3.22%  '
$ws.Range('E31').Value = '  -6.19%  '
$ws.Range('E32').Value = '  +4.11%  '
$ws.Range('E33').Value = '  -6.64%  '
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('E35').Value = '  -3.81%  '
$ws.Range('E36').Value = '  -3.01%  '
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('E38').Value = '  -5.55%  '
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('E41').Value = '  -11.31%  '
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('E43').Value = '  -10.68%  '
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('E45').Value = '  -6.82%  '
$ws.Range('E46').Value = '  -4.04%  '
$ws.Range('E47').Value = '  -6.09%  '
$ws.Range('E48').Value = '  -2.51%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('E51').Value = '  -0.55%  '
